# "Changed the nanogrid representation"
#
# 1) Re-cache the footer "datetimeFigureOut" date field (12/9/2019 -> 12/12/2019)
#    on the Slide Master and on every Custom Layout.
# 2) Recolor every nanogrid bar (Rectangle shape solid fill) from
#    212121 to 343434 on every slide, including rectangles nested inside
#    groups.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.Type -eq 14 -and $sh.PlaceholderFormat.Type -eq 16) {
            $sh.TextFrame.TextRange.Text = "12/12/2019"
        }
    }
}

# Update the date placeholder text on the slide master ...
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# ... and on every custom (slide) layout off that master.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

function Update-NanogridColor($shape) {
    if ($shape.Type -eq 6) {
        # msoGroup -> recurse into its GroupItems
        for ($k = 1; $k -le $shape.GroupItems.Count; $k++) {
            Update-NanogridColor $shape.GroupItems.Item($k)
        }
    } elseif ($shape.Fill.Type -eq 1 -and $shape.Fill.ForeColor.RGB -eq 0x212121) {
        $shape.Fill.ForeColor.RGB = 0x343434
    }
}

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        Update-NanogridColor $s.Shapes.Item($i)
    }
}
